$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new initiative rows were appended at the bottom of the sheet (rows 66
# and 67), both funded under axis "IV" with a "A criar" (to-be-created)
# project code, mirroring the formatting of the row immediately above them
# (row 64: s=11/3/3/6, wrapped text, row height 30).

# Row 66
$ws.Range("A64:D64").Copy()
$ws.Range("A66:D66").PasteSpecial(-4122)
$ws.Range("A66").Value = "A criar"
$ws.Range("B66").Value = "Melhoria da infraestrutura dos municípios – Fortalecimento do saneamento básico de Mário Campos"
$ws.Range("C66").Value = "IV"
$ws.Range("D66").Value = 50000000
$ws.Rows.Item(66).RowHeight = 30

# Row 67
$ws.Range("A64:D64").Copy()
$ws.Range("A67:D67").PasteSpecial(-4122)
$ws.Range("A67").Value = "A criar"
$ws.Range("B67").Value = "Elaboração de projetos rodoviários - Brumadinho-Mário Campos-BR381"
$ws.Range("C67").Value = "IV"
$ws.Range("D67").Value = 10000000
$ws.Rows.Item(67).RowHeight = 30

$excel.CutCopyMode = 0

$ws.Range("B72").Select()
